$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D, shifting D:K -> E:L
$ws.Columns("D").Insert()

# Copy number formats/styles from the (now shifted) old column D, which is now column E,
# onto the freshly inserted (blank) column D, for each of the three data blocks.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the latest reporting period figures.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 9352000
$ws.Range("D9").Value = 7355300
$ws.Range("D10").Value = 1996700
$ws.Range("D12").Value = 355200
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 12000
$ws.Range("D15").Value = 64700
$ws.Range("D17").Value = 8863000
$ws.Range("D18").Value = 489000
$ws.Range("D20").Value = -66800
$ws.Range("D21").Value = 712100
$ws.Range("D22").Value = 61900
$ws.Range("D23").Value = 360300
$ws.Range("D24").Value = 119400
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 240900
$ws.Range("D27").Value = 277000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 8500
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 66800
$ws.Range("D33").Value = 285500
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 285500
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 326100
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 880300
$ws.Range("D44").Value = 1908700
$ws.Range("D45").Value = 422300
$ws.Range("D46").Value = 3537400
$ws.Range("D47").Value = 400000
$ws.Range("D48").Value = 1373100
$ws.Range("D49").Value = 2068600
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 247300
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 7626400
$ws.Range("D57").Value = 865900
$ws.Range("D58").Value = 184200
$ws.Range("D59").Value = 1716600
$ws.Range("D60").Value = 2766700
$ws.Range("D61").Value = 1275300
$ws.Range("D62").Value = 590900
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 4693500
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 4477300
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 2932900
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 285500
$ws.Range("D83").Value = 289900
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 595900
$ws.Range("D91").Value = -203300
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -205500
$ws.Range("D96").Value = -47100
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -413300
$ws.Range("D101").Value = -18700
$ws.Range("D102").Value = -41600
